$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "{'int', 'any'}"
$ws.Range("E5").Value = "int"
$ws.Range("F5").Value = "Neutral"
$ws.Range("F5").Interior.Color = 42495
$ws.Range("E74").Value = "{'bool', 'empty'}"
$ws.Range("E75").Value = "bool"
$ws.Range("F75").Value = "Neutral"
$ws.Range("F75").Interior.Color = 42495
$ws.Range("E108").Value = "{'bool', 'any'}"
$ws.Range("E109").Value = "bool"
$ws.Range("F109").Value = "Neutral"
$ws.Range("F109").Interior.Color = 42495
$ws.Range("E138").Value = "{'Tuple[any]', 'Tuple[None]'}"
$ws.Range("E139").Value = "Tuple[any]"
$ws.Range("E252").Value = "{'bool', 'empty'}"
$ws.Range("E253").Value = "bool"
$ws.Range("F253").Value = "Neutral"
$ws.Range("F253").Interior.Color = 42495
$ws.Range("E268").Value = "{'bool', 'any'}"
$ws.Range("E269").Value = "bool"
$ws.Range("F269").Value = "Neutral"
$ws.Range("F269").Interior.Color = 42495
$ws.Range("E298").Value = "{'bool', 'any'}"
$ws.Range("E299").Value = "bool"
$ws.Range("F299").Value = "Neutral"
$ws.Range("F299").Interior.Color = 42495
$ws.Range("E371").Value = "{'Leaf', 'any'}"
$ws.Range("E372").Value = "Leaf"
$ws.Range("E580").Value = "{'TErr', 'Ok'}"
$ws.Range("E581").Value = "TErr"
$ws.Range("E582").Value = "{'Err', 'Ok', 'empty'}"
$ws.Range("E583").Value = "Err"
$ws.Range("E584").Value = "{'TErr', 'Ok'}"
$ws.Range("E585").Value = "TErr"
$ws.Range("E586").Value = "{'TErr', 'Ok'}"
$ws.Range("E587").Value = "TErr"
$ws.Range("E588").Value = "{'TErr', 'Ok'}"
$ws.Range("E589").Value = "TErr"
$ws.Range("E596").Value = "{'int', 'any', 'empty'}"
$ws.Range("E597").Value = "int"
$ws.Range("E602").Value = "{'bool', 'any'}"
$ws.Range("E603").Value = "bool"
$ws.Range("E604").Value = "{'bool', 'any'}"
$ws.Range("E605").Value = "bool"
$ws.Range("E612").Value = "{'int', 'empty'}"
$ws.Range("E613").Value = "int"
$ws.Range("F613").Value = "Neutral"
$ws.Range("F613").Interior.Color = 42495
$ws.Range("E624").Value = "{'TErr', 'Ok'}"
$ws.Range("E625").Value = "TErr"
$ws.Range("E628").Value = "{'TErr', 'Ok'}"
$ws.Range("E629").Value = "TErr"
$ws.Range("E632").Value = "{'TErr', 'Ok'}"
$ws.Range("E633").Value = "TErr"
$ws.Range("E634").Value = "{'TErr', 'Ok'}"
$ws.Range("E635").Value = "TErr"
$ws.Range("E638").Value = "{'TErr', 'Ok'}"
$ws.Range("E639").Value = "TErr"
$ws.Range("E640").Value = "{'TErr', 'Ok'}"
$ws.Range("E641").Value = "TErr"
$ws.Range("E652").Value = "{'TErr', 'Ok'}"
$ws.Range("E653").Value = "TErr"
$ws.Range("E881").Value = "{'Leaf', 'Node', 'any', 'empty'}"
$ws.Range("E882").Value = "Leaf"
$ws.Range("E969").Value = "{'bool', 'empty'}"
$ws.Range("E970").Value = "bool"
$ws.Range("F970").Value = "Neutral"
$ws.Range("F970").Interior.Color = 42495
$ws.Range("E971").Value = "{'bool', 'empty'}"
$ws.Range("E972").Value = "bool"
$ws.Range("F972").Value = "Neutral"
$ws.Range("F972").Interior.Color = 42495
$ws.Range("E975").Value = "{'bool', 'empty'}"
$ws.Range("E976").Value = "bool"
$ws.Range("F976").Value = "Neutral"
$ws.Range("F976").Interior.Color = 42495
$ws.Range("E1123").Value = "{'bool', 'empty'}"
$ws.Range("E1124").Value = "bool"
$ws.Range("F1124").Value = "Neutral"
$ws.Range("F1124").Interior.Color = 42495
$ws.Range("E1125").Value = "{'bool', 'empty'}"
$ws.Range("E1126").Value = "bool"
$ws.Range("E1127").Value = "{'bool', 'empty'}"
$ws.Range("E1128").Value = "bool"
$ws.Range("F1128").Value = "Neutral"
$ws.Range("F1128").Interior.Color = 42495
$ws.Range("E1129").Value = "{'bool', 'empty'}"
$ws.Range("E1130").Value = "bool"
$ws.Range("E1218").Value = "{'bool', 'empty'}"
$ws.Range("E1219").Value = "bool"
$ws.Range("F1219").Value = "Neutral"
$ws.Range("F1219").Interior.Color = 42495
$ws.Range("E1222").Value = "{'bool', 'empty'}"
$ws.Range("E1223").Value = "bool"
$ws.Range("F1223").Value = "Neutral"
$ws.Range("F1223").Interior.Color = 42495
$ws.Range("E1263").Value = "{'List[any]', 'empty'}"
$ws.Range("E1264").Value = "List[any]"
$ws.Range("F1264").Value = "Neutral"
$ws.Range("F1264").Interior.Color = 42495
$ws.Range("E1275").Value = "{'bool', 'empty'}"
$ws.Range("E1276").Value = "bool"
$ws.Range("F1276").Value = "Neutral"
$ws.Range("F1276").Interior.Color = 42495
$ws.Range("E1293").Value = "{'bool', 'empty'}"
$ws.Range("E1294").Value = "bool"
$ws.Range("F1294").Value = "Neutral"
$ws.Range("F1294").Interior.Color = 42495

# Row 1397: PyType Wins count changes from 299 to 282
$ws.Range("D1397").Value = 282

# Row 1398: "Scalpel Accuracy:" label moves from C1398 to E1398; value moves from D1398 to F1398 (366.56 -> 79.78)
$ws.Range("C1398").Value = ""
$ws.Range("D1398").Value = ""
$ws.Range("E1398").Value = "Scalpel Accuracy:"
$ws.Range("F1398").Value = 79.78

# Row 1399: text change and value change
$ws.Range("E1399").Value = "Accuracy vs PyType"
$ws.Range("F1399").Value = 36.17
